$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text type (many values look numeric,
# e.g. "10.00" or "0.0000250", and must not be auto-converted by Excel).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.200.13'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.798.66'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.96'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.31'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.23%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.80%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.452'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.49'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000250'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.79'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.431.43'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.787.97'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '68.152.04'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.58%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.46'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.115'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.09'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.57%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '461.26'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.72'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.700'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -3.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '82.98'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.75%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.57%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.947.21'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.65'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.23'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.36'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.02'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.78%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0998'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.53%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.30'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.80'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.989'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.63%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '43.66'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.68%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.41'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '151.95'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.36'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.58%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '396.58'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.36'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.51%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.81'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.66%  '
